$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "74.035.51"
$ws.Range("E2").Value = "  +7.49%  "

# Row 3
$ws.Range("D3").Value = "2.627.07"
$ws.Range("E3").Value = "  +7.63%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +14.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.203"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +19.31%  "

# Row 10
$ws.Range("D10").Value = "2.621.74"
$ws.Range("E10").Value = "  +7.44%  "

# Row 11
$ws.Range("E11").Value = "  +0.45%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.61%  "

# Row 15
$ws.Range("D15").Value = "3.106.38"
$ws.Range("E15").Value = "  +7.55%  "

# Row 16
$ws.Range("D16").Value = "73.689.02"
$ws.Range("E16").Value = "  +7.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.73%  "

# Row 18
$ws.Range("D18").Value = "2.636.50"
$ws.Range("E18").Value = "  +8.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +31.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.78%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.36%  "

# Row 28
$ws.Range("D28").Value = "2.757.34"
$ws.Range("E28").Value = "  +7.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.00%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0945"
$ws.Range("E30").Value = "  +15.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "525.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +22.17%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +19.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.34%  "

# Row 34
$ws.Range("E34").Value = "  +8.94%  "

# Row 35
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "

# Row 37
$ws.Range("E37").Value = "  +12.92%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.39%  "

# Row 39
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.94%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "161.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +24.05%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.39%  "

# Row 46
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.40%  "

# Row 47
$ws.Range("E47").Value = "  +3.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0853"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.38%  "

# Row 49
$ws.Range("E49").Value = "  +8.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.528"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +23.02%  "
